# Applies the commit: delete "ngach", insert "dian_tkt"
# (concretely: drop the stray _GoBack bookmark from the "ke hoach" paragraph,
#  re-flow the "Bo tri doan thanh tra gom" bullet's indent/spacing and add the
#  new "..., do <cb_cv> lam Truong doan" clause, and re-flow the "Thoi gian
#  tien hanh thanh tra" bullet's indent while moving the _GoBack bookmark
#  there.)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph "Can cu ke hoach phan bo thoi gian ... <nam_kh_tkt>;"
#    -> remove the _GoBack bookmark that sits right after the trailing ';'
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Paragraph "Bo tri doan thanh tra gom: <sl_cb> dong chi phong Thanh
#    tra - Kiem tra."
#    -> drop the explicit tab stop, bump the line spacing, switch to a
#       firstLine indent, give the <sl_cb> run (and the space after it)
#       the 14pt run size, and append ", do <cb_cv> lam Truong doan"
#       right before the final period.
# ---------------------------------------------------------------------
$pBoTri = $d.Paragraphs.Item(14)
$pBoTri.Format.TabStops.ClearAll()
$pBoTri.Format.LineSpacingRule = 5
$pBoTri.Format.LineSpacing = 13.8
$pBoTri.Format.LeftIndent = 0
$pBoTri.Format.FirstLineIndent = 28.35

$rSlCb = $pBoTri.Range.Duplicate()
$rSlCb.Find.Execute("<sl_cb> ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rSlCb.Font.Size = 14

$rTail = $pBoTri.Range.Duplicate()
$rTail.Find.Execute("Kiểm tra.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rTail.End - 1, $rTail.End - 1)
$insertPoint.InsertAfter(", do  <cb_cv> làm Trưởng đoàn")

# ---------------------------------------------------------------------
# 3) Paragraph "Thoi gian tien hanh thanh tra: nam <nam_ktra>."
#    -> drop the explicit tab stop, switch to a hanging indent, and move
#       the _GoBack bookmark to the start of this paragraph.
# ---------------------------------------------------------------------
$pThoiGian = $d.Paragraphs.Item(19)
$pThoiGian.Format.TabStops.ClearAll()
$pThoiGian.Format.LeftIndent = 35.45
$pThoiGian.Format.FirstLineIndent = -7.1

$goBackPoint = $d.Range($pThoiGian.Range.Start, $pThoiGian.Range.Start)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
